$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.252.15'
$ws.Range('E2').Value = '  +12.13%  '

$ws.Range('D3').Value = '1.823.70'
$ws.Range('E3').Value = '  +9.07%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.32%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.15'
$ws.Range('E5').Value = '  +4.41%  '

$ws.Range('E6').Value = '  +3.97%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.28%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.57'
$ws.Range('E8').Value = '  +6.89%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.51'
$ws.Range('E9').Value = '  +7.46%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.283'
$ws.Range('E10').Value = '  +7.18%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0669'
$ws.Range('E11').Value = '  +5.53%  '

$ws.Range('E12').Value = '  +2.64%  '

$ws.Range('D13').Value = '2.079.73'
$ws.Range('E13').Value = '  +8.68%  '

$ws.Range('D14').Value = '1.812.98'
$ws.Range('E14').Value = '  +8.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.645'
$ws.Range('E15').Value = '  +5.22%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '10.44'
$ws.Range('E16').Value = '  +3.00%  '

$ws.Range('D17').Value = '34.173.09'
$ws.Range('E17').Value = '  +11.83%  '

$ws.Range('E18').Value = '  +8.26%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.71'
$ws.Range('E19').Value = '  +5.47%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '259.10'
$ws.Range('E20').Value = '  +6.59%  '

$ws.Range('E21').Value = '  +4.23%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.23%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('E23').Value = '  +6.47%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.36'
$ws.Range('E24').Value = '  +2.73%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.22'
$ws.Range('E25').Value = '  +3.63%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.37'
$ws.Range('E26').Value = '  +0.36%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.66'
$ws.Range('E27').Value = '  +5.12%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.16'
$ws.Range('E28').Value = '  +7.50%  '

$ws.Range('E29').Value = '  +2.94%  '

$ws.Range('E30').Value = '  -0.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.88'
$ws.Range('E31').Value = '  +12.02%  '

$ws.Range('E32').Value = '  +4.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').Value = '  +5.62%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.55'
$ws.Range('E34').Value = '  +8.37%  '

$ws.Range('D35').Value = '1.548.87'
$ws.Range('E35').Value = '  +3.51%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.82'
$ws.Range('E36').Value = '  +3.88%  '

$ws.Range('E37').Value = '  +6.36%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '85.57'
$ws.Range('E38').Value = '  +2.32%  '

$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.628'
$ws.Range('E39').Value = '  +6.07%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0188'
$ws.Range('E40').Value = '  +5.70%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.82'
$ws.Range('E41').Value = '  +4.54%  '

$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.917'
$ws.Range('E42').Value = '  +9.66%  '

$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.34'
$ws.Range('E43').Value = '  +1.07%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.16'
$ws.Range('E44').Value = '  +9.46%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0524'
$ws.Range('E45').Value = '  +5.40%  '

$ws.Range('E46').Value = '  +4.44%  '

$ws.Range('D47').Value = '1.984.32'
$ws.Range('E47').Value = '  +9.76%  '

$ws.Range('B48').Value = 'MinaProtocolToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.06'
$ws.Range('E48').Value = '  +158.82%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.74'
$ws.Range('E49').Value = '  +3.39%  '

$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.05'
$ws.Range('E50').Value = '  +3.51%  '

$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.997'
$ws.Range('E51').Value = '  -0.34%  '
